# Final N-policy linking plots
# Adds two more iteration blocks (E:G = Iteration_1, H:J = Iteration_2)
# mirroring the existing Standalone block (B:D), and updates the numeric
# data rows 4-15 to their final, converged values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row 1: merged iteration labels -----------------------------
# Merge the new header regions first, *then* copy the style from an
# existing single (unmerged) header cell, A1, onto them. Doing the merge
# first avoids Excel re-deriving per-cell left/middle/right border
# variants for the merged block, so every header cell ends up sharing the
# same plain style index (bold, centered, bordered), just like B1:D1.
$ws.Range("E1:G1").Merge() | Out-Null
$ws.Range("H1:J1").Merge() | Out-Null

$ws.Range("A1").Copy() | Out-Null
$ws.Range("E1:J1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$ws.Range("E1").Value = "Iteration_1"
$ws.Range("H1").Value = "Iteration_2"

# --- Header row 2: year labels ------------------------------------------
$ws.Range("A2").Copy() | Out-Null
$ws.Range("E2:J2").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

# Force text storage (matching the existing 2030/2040/2050 labels in B2:D2,
# which are stored as text, not numbers) by applying a text number format
# before assigning the values.
$ws.Range("E2:J2").NumberFormat = "@"

$ws.Range("E2").Value = "2030"
$ws.Range("F2").Value = "2040"
$ws.Range("G2").Value = "2050"
$ws.Range("H2").Value = "2030"
$ws.Range("I2").Value = "2040"
$ws.Range("J2").Value = "2050"

# --- Data rows 4-15: final values for Standalone / Iteration_1 / Iteration_2 --
# Values are provided as strings and parsed to [double] because the script
# parser here does not accept scientific-notation numeric literals directly.
$values = @{
    4  = @("1184000.000000033", "1184000", "0", "1183999.999999957", "0", "0", "1183999.999999999", "0", "0")
    5  = @("0", "0", "0", "0", "0", "0", "0", "0", "0")
    6  = @("7.412782555161258e-09", "0", "0", "0", "1179586.812507884", "1179258.282746874", "0", "1179593.288345533", "1179267.814186228")
    7  = @("9.761630826532434e-11", "554.2455973799724", "1177607.500706655", "0", "0", "0", "0", "0", "0")
    8  = @("0", "0", "0", "0", "0", "0", "0", "0", "0")
    9  = @("0", "0", "0", "0", "0", "0", "0", "0", "0")
    10 = @("0", "0", "0", "0", "0", "0", "0", "0", "0")
    11 = @("0", "0", "0", "0", "0", "0", "0", "0", "0")
    12 = @("0", "0", "0", "0", "0", "0", "0", "0", "0")
    13 = @("0", "0", "0", "0", "0", "0", "0", "0", "0")
    14 = @("0", "0", "0", "0", "0", "0", "0", "0", "0")
    15 = @("0", "0", "0", "0", "0", "0", "0", "0", "0")
}

foreach ($r in $values.Keys) {
    $rowVals = $values[$r]
    for ($i = 0; $i -lt $rowVals.Length; $i++) {
        $col = 2 + $i   # B=2 .. J=10
        $ws.Cells.Item($r, $col).Value = [double]$rowVals[$i]
    }
}
